$wb = $excel.ActiveWorkbook

# Rename the first sheet from "Sheet1" to "TestEnv"
$wsTest = $wb.Worksheets.Item("Sheet1")
$wsTest.Name = "TestEnv"

# Make the TestEnv sheet the active/selected tab instead of PreProdEnv
$wsTest.Activate()
